# Applies the changes described by the target diff:
#  - moves the "active tab" / selected sheet from Input to Transactions
#  - updates the active-cell selection on every sheet
#  - updates a handful of computed repayment-schedule / summary values
#  - removes a few now-blank helper cells on the Repayment Schedule sheet

$wb = $excel.ActiveWorkbook

$wsInput      = $wb.Worksheets.Item("Input")
$wsSummary    = $wb.Worksheets.Item("Summary")
$wsSchedule   = $wb.Worksheets.Item("Repayment Schedule")
$wsTxns       = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------------
# Summary sheet: correct a couple of recomputed currency figures
# ---------------------------------------------------------------------------
$wsSummary.Range("A3").Value = 672.36
$wsSummary.Range("E3").Value = 570.44000000000005

# ---------------------------------------------------------------------------
# Repayment Schedule sheet: recomputed figures across several rows
# ---------------------------------------------------------------------------
$wsSchedule.Range("F6").Value  = 812.55
$wsSchedule.Range("G6").Value  = 6808.56
$wsSchedule.Range("H6").Value  = 75.17

$wsSchedule.Range("F7").Value  = 818.33
$wsSchedule.Range("G7").Value  = 5990.23
$wsSchedule.Range("H7").Value  = 69.39

$wsSchedule.Range("F8").Value  = 828.64
$wsSchedule.Range("G8").Value  = 5161.59
$wsSchedule.Range("H8").Value  = 59.08

$wsSchedule.Range("G9").Value  = 4326.4799999999996

$wsSchedule.Range("F10").Value = 843.63
$wsSchedule.Range("G10").Value = 3482.85
$wsSchedule.Range("H10").Value = 44.09

$wsSchedule.Range("F11").Value = 853.37
$wsSchedule.Range("G11").Value = 2629.48
$wsSchedule.Range("H11").Value = 34.35

$wsSchedule.Range("G12").Value = 1768.56

$wsSchedule.Range("F13").Value = 870.28
$wsSchedule.Range("G13").Value = 898.28
$wsSchedule.Range("H13").Value = 17.440000000000001

$wsSchedule.Range("F14").Value = 898.28
$wsSchedule.Range("K14").Value = 907.44
$wsSchedule.Range("P14").Value = 907.44

# Fully clear the now-empty helper cells in row 2 so they disappear from the
# sheet (rather than remaining as blank, styled cells).
$wsSchedule.Range("A2").Clear()
$wsSchedule.Range("E2").Clear()
$wsSchedule.Range("N2").Clear()
$wsSchedule.Range("O2").Clear()

# ---------------------------------------------------------------------------
# Sheet view / selection updates, finishing with Transactions as the active
# (selected) sheet/tab, matching the new workbookView activeTab.
# ---------------------------------------------------------------------------
$wsInput.Activate()
$wsInput.Range("C19").Select()

$wsSummary.Activate()
$wsSummary.Range("E29").Select()

$wsSchedule.Activate()
$wsSchedule.Range("I10").Select()

$wsTxns.Activate()
$wsTxns.Range("B3").Select()
